{"js": "const replacements = [\n  ['2024-01-09 Tuesday', '2024-01-10 Wednesday'],\n  ['955\u00d75=', '425\u00d77='],\n  ['451\u00d79=', '213\u00d77='],\n  ['121\u00d76=', '555\u00d75='],\n  ['566\u00d73=', '956\u00d72='],\n  ['375\u00d74=', '940\u00d75='],\n  ['653\u00d73=', '766\u00d74='],\n  ['650\u00d74=', '108\u00d77='],\n  ['927\u00d78=', '286\u00d78='],\n  ['231\u00d73=', '148\u00d79='],\n  ['496\u00d79=', '651\u00d72='],\n  ['544\u00d72=', '504\u00d77='],\n  ['857\u00d78=', '472\u00d72='],\n  ['864\u00d72=', '988\u00d73='],\n  ['390\u00d79=', '131\u00d75='],\n  ['743\u00d74=', '521\u00d76='],\n  ['613\u00d75=', '798\u00d76='],\n  ['854\u00d79=', '988\u00d76='],\n  ['182\u00d76=', '625\u00d73='],\n  ['383\u00d74=', '186\u00d76='],\n  ['980\u00d76=', '404\u00d74='],\n  ['978\u00d73=', '780\u00d77='],\n  ['993\u00d76=', '903\u00d79='],\n  ['733\u00d74=', '838\u00d79='],\n  ['541\u00d76=', '259\u00d73='],\n  ['400\u00d78=', '287\u00d74='],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(replace, 'Replace');\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$replacements = @(\n    @('2024-01-09 Tuesday', '2024-01-10 Wednesday'),\n    @('955\u00d75=', '425\u00d77='),\n    @('451\u00d79=', '213\u00d77='),\n    @('121\u00d76=', '555\u00d75='),\n    @('566\u00d73=', '956\u00d72='),\n    @('375\u00d74=', '940\u00d75='),\n    @('653\u00d73=', '766\u00d74='),\n    @('650\u00d74=', '108\u00d77='),\n    @('927\u00d78=', '286\u00d78='),\n    @('231\u00d73=', '148\u00d79='),\n    @('496\u00d79=', '651\u00d72='),\n    @('544\u00d72=', '504\u00d77='),\n    @('857\u00d78=', '472\u00d72='),\n    @('864\u00d72=', '988\u00d73='),\n    @('390\u00d79=', '131\u00d75='),\n    @('743\u00d74=', '521\u00d76='),\n    @('613\u00d75=', '798\u00d76='),\n    @('854\u00d79=', '988\u00d76='),\n    @('182\u00d76=', '625\u00d73='),\n    @('383\u00d74=', '186\u00d76='),\n    @('980\u00d76=', '404\u00d74='),\n    @('978\u00d73=', '780\u00d77='),\n    @('993\u00d76=', '903\u00d79='),\n    @('733\u00d74=', '838\u00d79='),\n    @('541\u00d76=', '259\u00d73='),\n    @('400\u00d78=', '287\u00d74='),\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 1\n    $range.Find.MatchCase = $true\n    $range.Find.MatchWholeWord = $false\n    $range.Find.MatchWildcards = $false\n    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
